$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 44 (this shifts the existing rows 44:191 down
# to 46:193, growing the table from 191 to 193 rows).
$ws.Rows("44:45").Insert()

# The two newly-inserted blank rows (now 44:45) become a new weekly entry.
# Its Primera/Segunda price data is identical to the entry that is now
# sitting right below it (rows 46:47, which is what used to be rows 44:45
# before the insert), so copy that data up into the new rows.
$ws.Range("A46:R47").Copy($ws.Range("A44:R45"))

# ...except the new entry carries a new date.
$ws.Range("D44").Value = 44453
$ws.Range("D45").Value = 44453
